$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New payment row appended below the header (row 1), matching the
# existing columns: phone, original_amount, birthday_discount,
# reward_discount, points_redeemed, final_amount, method, timestamp
$ws.Range("A2").Value = "'76442711"
$ws.Range("B2").Value = 408
$ws.Range("C2").Value = 61.2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 346.8
$ws.Range("G2").Value = "Cash"
$ws.Range("H2").Value = "'2025-08-20T08:52:14"

# Strip the quote-prefix "Normal"-but-text style Excel would otherwise
# assign to the force-text cells so the new row stays unstyled, like the
# rest of the data rows.
$ws.Range("A2").Style = "Normal"
$ws.Range("H2").Style = "Normal"
